$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.254.43"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.906.05"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5234"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.94%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3782"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07258"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9007"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08210"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.55%  "
$ws.Range("D13").Value = "1.910.47"
$ws.Range("E13").Value = "  +2.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "95.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.288"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.0000"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008595"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("D20").Value = "27.298.60"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.066"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "2.162.86"
$ws.Range("E22").Value = "  +2.15%  "
$ws.Range("E23").Value = "  +3.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.465"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.303"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.746"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.985"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.811"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09214"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8054"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.40%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05073"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.242"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.956"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.333"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.566"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.44%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5724"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01982"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.064"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.627"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4838"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.618"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.23%  "
